$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C6").Value = 8623
$ws.Range("C7:C12").Value = 8004
$ws.Range("C13:C252").Value = 7586
